# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# across all 8 item-sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""

$ws.Range("H11").Value = 171.44444
$ws.Range("I11").Value = 171.44444
$ws.Range("K11").Value = 171.44444
$ws.Range("M11").Value = -31.44443999999999

$ws.Range("H17").Value = 3501.611
$ws.Range("J17").Value = 2868.9375
$ws.Range("L17").Value = 8606.8125
$ws.Range("N17").Value = -8942.8125

$ws.Range("H18").Value = 10383.429
$ws.Range("I18").Value = 5875.231
$ws.Range("J18").Value = 17709.25
$ws.Range("K18").Value = 5875.231
$ws.Range("L18").Value = 17709.25
$ws.Range("M18").Value = -5591.231
$ws.Range("N18").Value = -18277.25

$ws.Range("H51").Value = 6000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H74").Value = 3749.6667
$ws.Range("I74").Value = 2499.6667
$ws.Range("K74").Value = 2499.6667
$ws.Range("M74").Value = -1563.6667

$ws.Range("H77").Value = 3749.6667
$ws.Range("I77").Value = 2499.6667
$ws.Range("K77").Value = 12498.3335
$ws.Range("M77").Value = -7818.333500000001

$ws.Range("H92").Value = 261.41177
$ws.Range("I92").Value = 215.375
$ws.Range("J92").Value = 998
$ws.Range("K92").Value = 215.375
$ws.Range("L92").Value = 998
$ws.Range("M92").Value = 1032.625
$ws.Range("N92").Value = -3494

$ws.Range("H116").Value = 23320
$ws.Range("I116").Value = 50750
$ws.Range("J116").Value = 5033.3335
$ws.Range("K116").Value = 50750
$ws.Range("L116").Value = 5033.3335
$ws.Range("M116").Value = -47308
$ws.Range("N116").Value = -11917.3335

$ws.Range("H129").Value = 879.8036
$ws.Range("J129").Value = 899.549
$ws.Range("L129").Value = 2698.647
$ws.Range("N129").Value = -12698.647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6171.735
$ws.Range("I32").Value = 4028
$ws.Range("K32").Value = 4028
$ws.Range("M32").Value = -3741

$ws.Range("H61").Value = 5790.231
$ws.Range("I61").Value = 6376.8237
$ws.Range("K61").Value = 6376.8237
$ws.Range("M61").Value = -6164.8237

$ws.Range("H63").Value = 1596.3334
$ws.Range("I63").Value = 1535.6
$ws.Range("K63").Value = 1535.6
$ws.Range("M63").Value = -849.5999999999999

$ws.Range("H66").Value = 1596.3334
$ws.Range("I66").Value = 1535.6
$ws.Range("K66").Value = 7678
$ws.Range("M66").Value = -4246

$ws.Range("H74").Value = 3113.8572
$ws.Range("J74").Value = 3959.4
$ws.Range("L74").Value = 3959.4
$ws.Range("N74").Value = -5707.4

$ws.Range("H77").Value = 3113.8572
$ws.Range("J77").Value = 3959.4
$ws.Range("L77").Value = 19797
$ws.Range("N77").Value = -28533

$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100722

$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102496

$ws.Range("H122").Value = 1150.4286
$ws.Range("I122").Value = 861.46155
$ws.Range("J122").Value = 1620
$ws.Range("K122").Value = 2584.38465
$ws.Range("L122").Value = 4860
$ws.Range("M122").Value = -134.38465
$ws.Range("N122").Value = -9760

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

$ws.Range("H132").Value = 1670.258
$ws.Range("I132").Value = 1325.0869
$ws.Range("J132").Value = 2662.625
$ws.Range("K132").Value = 3975.2607
$ws.Range("L132").Value = 7987.875
$ws.Range("M132").Value = -1445.2607
$ws.Range("N132").Value = -13047.875

$ws.Range("H136").Value = 5790.231
$ws.Range("I136").Value = 6376.8237
$ws.Range("K136").Value = 19130.4711
$ws.Range("M136").Value = -16580.4711

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26457.125
$ws.Range("I82").Value = 10052.333
$ws.Range("J82").Value = 36300
$ws.Range("K82").Value = 10052.333
$ws.Range("L82").Value = 36300
$ws.Range("M82").Value = -9669.333000000001
$ws.Range("N82").Value = -37066

$ws.Range("H85").Value = 26457.125
$ws.Range("I85").Value = 10052.333
$ws.Range("J85").Value = 36300
$ws.Range("K85").Value = 10052.333
$ws.Range("L85").Value = 36300
$ws.Range("M85").Value = -8726.333000000001
$ws.Range("N85").Value = -38952

$ws.Range("H86").Value = 339218.34
$ws.Range("I86").Value = 13900
$ws.Range("J86").Value = 501877.5
$ws.Range("K86").Value = 13900
$ws.Range("L86").Value = 501877.5
$ws.Range("M86").Value = -12777
$ws.Range("N86").Value = -504123.5

$ws.Range("H89").Value = 339218.34
$ws.Range("I89").Value = 13900
$ws.Range("J89").Value = 501877.5
$ws.Range("K89").Value = 69500
$ws.Range("L89").Value = 2509387.5
$ws.Range("M89").Value = -63884
$ws.Range("N89").Value = -2520619.5

$ws.Range("H105").Value = 2286.7778
$ws.Range("I105").Value = 2015.7826
$ws.Range("K105").Value = 2015.7826
$ws.Range("M105").Value = -268.7826

$ws.Range("H107").Value = 904.9375
$ws.Range("I107").Value = 632.8333
$ws.Range("J107").Value = 1068.2
$ws.Range("K107").Value = 632.8333
$ws.Range("L107").Value = 1068.2
$ws.Range("M107").Value = 1287.1667
$ws.Range("N107").Value = -4908.2

$ws.Range("H134").Value = 4729.9697
$ws.Range("I134").Value = 5252.731
$ws.Range("K134").Value = 15758.193
$ws.Range("M134").Value = -13223.193

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9166.666999999999

$ws.Range("H122").Value = 3305.6667
$ws.Range("I122").Value = 2357.3
$ws.Range("J122").Value = 5202.4
$ws.Range("K122").Value = 7071.900000000001
$ws.Range("L122").Value = 15607.2
$ws.Range("M122").Value = -4621.900000000001
$ws.Range("N122").Value = -20507.2

$ws.Range("H134").Value = 3538.0908
$ws.Range("I134").Value = 3000.6667
$ws.Range("J134").Value = 5956.5
$ws.Range("K134").Value = 9002.000100000001
$ws.Range("L134").Value = 17869.5
$ws.Range("M134").Value = -6467.000100000001
$ws.Range("N134").Value = -22939.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 552.6667
$ws.Range("J34").Value = 1104.4
$ws.Range("L34").Value = 3313.2
$ws.Range("N34").Value = -3481.2

$ws.Range("H38").Value = 509.66666
$ws.Range("I38").Value = 139.375
$ws.Range("K38").Value = 418.125
$ws.Range("M38").Value = -71.125

$ws.Range("H68").Value = 1101
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1101
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3303
$ws.Range("M68").Value = ""
$ws.Range("N68").Value = -4925

$ws.Range("H71").Value = 1101
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1101
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 9909
$ws.Range("M71").Value = ""
$ws.Range("N71").Value = -18021

$ws.Range("H113").Value = 5661.048
$ws.Range("J113").Value = 910.05554
$ws.Range("L113").Value = 2730.16662
$ws.Range("N113").Value = -7070.16662

$ws.Range("H129").Value = 81283
$ws.Range("J129").Value = 121469.836
$ws.Range("L129").Value = 364409.508
$ws.Range("N129").Value = -374409.508

$ws.Range("H131").Value = 16655.342
$ws.Range("J131").Value = 17412.023
$ws.Range("L131").Value = 52236.069
$ws.Range("N131").Value = -62316.069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23800
$ws.Range("I46").Value = 20000
$ws.Range("K46").Value = 20000
$ws.Range("M46").Value = -19844

$ws.Range("H97").Value = 1494
$ws.Range("I97").Value = 1680
$ws.Range("K97").Value = 1680
$ws.Range("M97").Value = -1184

$ws.Range("H122").Value = 1544.6842
$ws.Range("I122").Value = 1597.8334
$ws.Range("J122").Value = 1453.5714
$ws.Range("K122").Value = 4793.5002
$ws.Range("L122").Value = 4360.7142
$ws.Range("M122").Value = -2343.5002
$ws.Range("N122").Value = -9260.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 479166.66
$ws.Range("J2").Value = 250000
$ws.Range("L2").Value = 250000
$ws.Range("N2").Value = -250224

$ws.Range("H22").Value = 2651.2856
$ws.Range("I22").Value = 601
$ws.Range("J22").Value = 2993
$ws.Range("K22").Value = 601
$ws.Range("L22").Value = 2993
$ws.Range("N22").Value = -3583
$ws.Range("M22").Value = -306

$ws.Range("H27").Value = 2651.2856
$ws.Range("I27").Value = 601
$ws.Range("J27").Value = 2993
$ws.Range("K27").Value = 601
$ws.Range("L27").Value = 2993
$ws.Range("N27").Value = -3207
$ws.Range("M27").Value = -494

$ws.Range("H136").Value = 7400
$ws.Range("I136").Value = 5666.6665
$ws.Range("K136").Value = 16999.9995
$ws.Range("M136").Value = -14449.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 49980
$ws.Range("I56").Value = 49980
$ws.Range("K56").Value = 49980
$ws.Range("M56").Value = -49266

$ws.Range("H126").Value = 1860.5714
$ws.Range("I126").Value = 1806.1111
$ws.Range("K126").Value = 5418.3333
$ws.Range("M126").Value = -2948.3333

$ws.Range("H136").Value = 24156072
$ws.Range("I136").Value = 34723350
$ws.Range("J136").Value = 2292.7144
$ws.Range("K136").Value = 104170050
$ws.Range("L136").Value = 6878.1432
$ws.Range("M136").Value = -104167500
$ws.Range("N136").Value = -11978.1432

$ws.Range("H138").Value = 620000
$ws.Range("J138").Value = 620000
$ws.Range("L138").Value = 620000
$ws.Range("N138").Value = -630280

Write-Host "Applied market-data refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
